$wb = $excel.ActiveWorkbook

# --- "Prix Spot" sheet: a new day (05-dec) was inserted as a column right
#     before the existing "01-oct." column (EH), shifting everything from
#     EH onward one column to the right (EH->EI ... FL->FM). ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$insertCol = $wsPrix.Range("EH1:EH25").EntireColumn
$insertCol.Insert()

# Header for the freshly inserted column.
$wsPrix.Range("EH1").Value = "05-dec"

# The new day has no data yet for any of the 24 hourly rows.
$wsPrix.Range("EH2:EH25").Value = "-"

# --- "Gaz" sheet: append the new day's price as row 168. ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A168").Value = "'2025-12-03"
$wsGaz.Range("B168").Value = 27.05

# --- "CO2" sheet: append the new day's price as row 168. ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A168").Value = "'2025-12-03"
$wsCO2.Range("B168").Value = 81.35
